# Generate Report for handoff
# The file "7962ae02-cf9c-42a6-8dfd-ad1f7df73638.md" (row 3 on every sheet) has been
# handed off again: its status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and the per-locale "Latest Handoff Datetime" is refreshed to
# reflect the new handoff generation time.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: File Name / zh-cn / de-de summary columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: Status + Latest Handoff Datetime for the same file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("D3").Value = "2016-01-26 05:29:45"

# --- de-de sheet: Status + Latest Handoff Datetime for the same file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("D3").Value = "2016-01-26 05:29:57"
